$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 13 ("Start testing" section header),
# shifting everything below it down by one row.
$ws.Rows.Item(13).Insert()

# Fill in the new instruction row: a reminder + the literal <script> snippet
# that needs to be added to MOLGENIS' header.
$ws.Range("A13").Value() = "Make sure d3.js is in the header of MOLGENIS, above jQuery!"
$ws.Range("B13").Value() = '<script src="https://cdnjs.cloudflare.com/ajax/libs/d3/3.5.6/d3.min.js" charset="utf-8"></script>'

# Match the row height used for other note rows, and keep B13 on the default
# (unstyled) look rather than the bold/section style picked up from the Insert.
$ws.Range("A13:B13").RowHeight = 15
$ws.Range("B13").Style = "Normal"

# Move the active selection to the newly added row (matches how the file
# is left selected after the edit, instead of the stale far-down selection).
$null = $ws.Range("B15").Select()
